# The deck has a single table (slide 6, shape 2 - "Sources of finance").
# Re-apply a different (built-in) table style to it, swapping it from the
# custom "Table_0" style to the built-in Medium-Style-2/Accent-1 style.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)
$shp = $s.Shapes.Item(2)

if ($shp.HasTable) {
    $tbl = $shp.Table
    $tbl.ApplyStyle("{4877FA82-A40D-4417-8CC3-4D921F6B7F2D}")
}
